# Update "想去人数" (want-to-go count) figures that changed between
# the previous and current data pull (gh-pages output regenerated at 456a3b4).
#
# Sheet "展览"    (sheet1): F4 411->412, F5 452->453, F7 2478->2479, F9 6532->6541
# Sheet "全部类型" (sheet4): F4 411->412, F5 452->453, F9 2478->2479, F11 6532->6541

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"    = @{ "F4" = 412; "F5" = 453; "F7" = 2479; "F9" = 6541 }
    "全部类型" = @{ "F4" = 412; "F5" = 453; "F9" = 2479; "F11" = 6541 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
